# Swap the contents of row 15 and row 16 for the columns that differ
# between them (A, B, E, F, G, H, K, L, M, N, Q, R, AC). Columns that are
# identical between the two rows (C, D, P, S, T, U, V, W, Y, AA, AD, AE,
# AG, AT, AW, AX, AY, ...) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 15
$row2 = 16

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $val1 = $ws.Range($addr1).Value2
    $val2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $val2
    $ws.Range($addr2).Value2 = $val1
}

# K15:N15 are empty inline-string cells in the source row (row 15) that
# move to row 16 after the swap; row 15 loses them (becomes blank/empty).
foreach ($col in @("K", "L", "M", "N")) {
    $ws.Range("$col$row1").Value2 = $null
    $ws.Range("$col$row2").Value2 = ""
}

# AC15 "ringhack äldre" moves to AC16; AC15 becomes empty.
$acVal = $ws.Range("AC$row1").Value2
$ws.Range("AC$row1").Value2 = $null
$ws.Range("AC$row2").Value2 = $acVal
